$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solar Park")

# Clear the user-entered "Value" column data (rows 2-38). The header/label
# column (B) and special computed rows (39-47) are left untouched; clearing
# these cells also causes the sharedStrings table to shrink/renumber and the
# dependent formulas in C42:C47 to recompute (blank / #VALUE! / #N/A) on
# their own once their source cells are emptied.
$ws.Range("C2:C38").ClearContents()

# Rows 5, 17, 22 and 31 had manually "grown" heights (wrapped multi-line
# values). Once the text is cleared there is nothing to wrap, so restore the
# default row height via AutoFit instead of leaving the stale explicit
# height behind.
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(31).AutoFit()

# Rows 41 and 48 were hidden helper rows; unhide them.
$ws.Rows.Item(41).Hidden = $false
$ws.Rows.Item(48).Hidden = $false

# Reset the view: scroll back to the top (drop the stale topLeftCell) and
# move the selection to C6.
$ws.Activate()
$ws.Range("C6").Select()
